$d = $word.ActiveDocument

$replacements = @(
    @{old="89÷2="; new="96÷3="},
    @{old="98÷3="; new="62÷9="},
    @{old="82÷9="; new="83÷3="},
    @{old="52÷2="; new="99÷3="},
    @{old="91÷3="; new="76÷3="},
    @{old="92÷4="; new="57÷5="},
    @{old="46÷2="; new="88÷5="},
    @{old="92÷7="; new="92÷5="},
    @{old="41÷4="; new="46÷8="},
    @{old="65÷2="; new="54÷2="},
    @{old="55÷5="; new="33÷4="},
    @{old="73÷7="; new="76÷4="},
    @{old="53÷7="; new="34÷7="},
    @{old="47÷7="; new="98÷5="},
    @{old="59÷2="; new="23÷9="},
    @{old="22÷8="; new="20÷8="},
    @{old="34÷4="; new="18÷7="},
    @{old="53÷9="; new="33÷5="},
    @{old="37÷6="; new="98÷8="},
    @{old="57÷8="; new="34÷3="},
    @{old="68÷4="; new="11÷7="},
    @{old="39÷7="; new="19÷2="},
    @{old="80÷3="; new="87÷7="},
    @{old="48÷4="; new="58÷8="},
    @{old="98÷9="; new="37÷2="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
